# The sheet is a long, flat data table (one row per market observation) for
# "Haba" at "Femacal de La Calera". This edit inserts one new observation
# row right above the current row 127, pushing every row below it down by
# one (old row 127 becomes row 128, ..., old row 204 becomes row 205).
#
# The new row 127 reuses most of the surrounding context (mercado, region,
# codreg, categoria id/categoria, variedad, calidad, unidad de
# comercializacion's "origen", kg/unidades, clasificacion) and only differs
# in fecha, volumen, precio minimo/maximo/promedio, unidad de
# comercializacion and precio $/kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 127:204 down to 128:205, opening up a blank row 127.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(127, 1).Value = 3
$ws.Cells.Item(127, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(127, 3).Value = "Coquimbo"
$ws.Cells.Item(127, 4).Value = 44830
$ws.Cells.Item(127, 5).Value = 5
$ws.Cells.Item(127, 6).Value = 100112026
$ws.Cells.Item(127, 7).Value = "Haba"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 100
$ws.Cells.Item(127, 11).Value = 10000
$ws.Cells.Item(127, 12).Value = 11000
$ws.Cells.Item(127, 13).Value = 10450
$ws.Cells.Item(127, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(127, 15).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(127, 16).Value = 418
$ws.Cells.Item(127, 17).Value = 25
$ws.Cells.Item(127, 18).Value = "Hortaliza"
